$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $text) {
    $cell = $t.Cell($row, $col)
    $r = $cell.Range
    # Trim trailing cell-mark / paragraph-mark characters Word appends to cell.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

Set-CellText 1 1 "149×6=894"
Set-CellText 1 2 "362×3=1086"
Set-CellText 1 3 "561×7=3927"
Set-CellText 1 4 "387×2=774"
Set-CellText 1 5 "678×9=6102"

Set-CellText 5 1 "942×7=6594"
Set-CellText 5 2 "432×6=2592"
Set-CellText 5 3 "934×9=8406"
Set-CellText 5 4 "251×5=1255"
Set-CellText 5 5 "134×6=804"

Set-CellText 10 1 "695×8=5560"
Set-CellText 10 2 "842×7=5894"
Set-CellText 10 3 "230×2=460"
Set-CellText 10 4 "264×7=1848"
Set-CellText 10 5 "243×8=1944"

Set-CellText 15 1 "144×2=288"
Set-CellText 15 2 "625×6=3750"
Set-CellText 15 3 "784×6=4704"
Set-CellText 15 4 "157×5=785"
Set-CellText 15 5 "408×4=1632"

Set-CellText 20 1 "406×8=3248"
Set-CellText 20 2 "902×5=4510"
Set-CellText 20 3 "517×3=1551"
Set-CellText 20 4 "873×2=1746"
Set-CellText 20 5 "387×6=2322"
